# Updates the cryptos list (price / 1h volume columns, plus a couple of
# coin-name/link cell swaps) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "@" (Text) number format before writing numeric-looking price
# strings so Excel keeps them as text (preserving formats like trailing
# zeros, "0.500", "120.00", etc.) instead of silently converting them to
# numbers.
$ws.Range("D2").Value = "65.001.03"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "3.139.58"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.56"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.42"
$ws.Range("E6").Value = "  +5.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.132.92"
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.01"
$ws.Range("E10").Value = "  +14.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.00"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000223"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "3.646.98"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").Value = "64.965.81"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "529.32"
$ws.Range("E18").Value = "  +10.26%  "
$ws.Range("D19").Value = "3.140.49"
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.75"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("E21").Value = "  +2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").Value = "  +5.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("E23").Value = "  +4.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.75"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.67"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.69"
$ws.Range("E27").Value = "  +14.58%  "
$ws.Range("E28").Value = "  +3.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.18"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "559.30"
$ws.Range("E34").Value = "  +12.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0444"
$ws.Range("E37").Value = "  +9.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.77"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0817"
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("E40").Value = "  +12.72%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.068.84"
$ws.Range("E41").Value = "  +7.33%  "
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.27"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.257"
$ws.Range("E44").Value = "  +7.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  +7.89%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.03"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "120.00"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "0.0₃0525"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.10"
$ws.Range("E51").Value = "  +3.87%  "
